# Apply "papers and update databses" changes:
# - Add column F (PDF file names) for rows 3-19
# - Update selected cell in sheet view from D18 to D10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("publications")

# Row -> PDF filename mapping (column F), rows 3..19
$pdfMap = @{
    3  = "2020_TCDS_conflict.pdf"
    4  = "2019_TVCG_VR Fall.pdf"
    5  = "2019_VRCAI_sickness.pdf"
    6  = "2019_CHI_HapticImmersion.pdf"
    7  = "2018_ACCESS_PredictionError.pdf"
    8  = "2017_CHI_TrussFab.pdf"
    9  = "2016_UIST_Metamaterial.pdf"
    10 = "2016_PG_icon.pdf"
    11 = "2016_I3D_history.pdf"
    12 = "2015_UIST_LaserStacker.pdf"
    13 = "2015_UIST_Protopiper.pdf"
    14 = "2015_CHI_platener.pdf"
    15 = "2014_SA_Autocomplete.pdf"
    16 = "2014_CHI_viewpoint.pdf"
    17 = "2013_VRST_Faceton.pdf"
    18 = "2013_PG_Splattering.pdf"
    19 = "2011_SIG_revision.pdf"
}

# Insert in the same order the original author typed the values (bottom rows
# first, with rows 11/12 swapped) so that the shared-string table gets built
# up in the exact sequence recorded in the target workbook.
$rowOrder = @(19, 18, 17, 16, 15, 14, 13, 11, 12, 10, 9, 8, 7, 6, 5, 4, 3)
foreach ($row in $rowOrder) {
    $ws.Range("F$row").Value = $pdfMap[$row]
}

# Update the active selection in the sheet view
$ws.Range("D10").Select()
